$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.1272849253844004
$ws.Range("D3").Value = 0.1272908219143504
$ws.Range("D4").Value = 0.1796008640945683
$ws.Range("D5").Value = 0.1722192543128125
$ws.Range("D6").Value = 0.1389108337282178
$ws.Range("D7").Value = 0.1272876666760212
$ws.Range("D8").Value = 0.1274056338896295
